# Applies scraped-data update for Linea 141 (LP1912 / LP1912-215 / 6203-6173)
# commit: Horarios actualizados Linea 141 - 182
# Updates "Ultima actualizacion" / "Total filas" headers and refreshes the
# detail rows (sheet1 "LP1912" and sheet3 "6203-6173") with the latest scrape,
# including new rows appended at the end of each sheet's table, plus the
# header timestamp refresh on sheet2 "LP1912-215" (no row-level changes there).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 20:12:03'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 356'
$ws1.Cells.Item(137, 1).Value = '11:52:01'
$ws1.Cells.Item(137, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(137, 4).Value = 45
$ws1.Cells.Item(138, 1).Value = '11:47:17'
$ws1.Cells.Item(138, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(138, 4).Value = 50
$ws1.Cells.Item(160, 1).Value = '11:34:59'
$ws1.Cells.Item(160, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(160, 4).Value = 118
$ws1.Cells.Item(161, 1).Value = '12:11:52'
$ws1.Cells.Item(161, 3).Value = '14_ABASTO'
$ws1.Cells.Item(161, 4).Value = 81
$ws1.Cells.Item(216, 1).Value = '15:17:33'
$ws1.Cells.Item(216, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(216, 4).Value = 48
$ws1.Cells.Item(217, 1).Value = '14:12:26'
$ws1.Cells.Item(217, 3).Value = '14_ABASTO'
$ws1.Cells.Item(217, 4).Value = 113
$ws1.Cells.Item(259, 1).Value = '16:44:12'
$ws1.Cells.Item(259, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(259, 4).Value = 64
$ws1.Cells.Item(261, 1).Value = '17:48:33'
$ws1.Cells.Item(261, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(261, 4).Value = 0
$ws1.Cells.Item(283, 1).Value = '16:37:06'
$ws1.Cells.Item(283, 3).Value = '15X38_ABASTO'
$ws1.Cells.Item(283, 4).Value = 119
$ws1.Cells.Item(284, 1).Value = '17:13:39'
$ws1.Cells.Item(284, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(284, 4).Value = 83
$ws1.Cells.Item(305, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(306, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(327, 1).Value = '18:52:02'
$ws1.Cells.Item(327, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(327, 4).Value = 75
$ws1.Cells.Item(328, 1).Value = '19:48:11'
$ws1.Cells.Item(328, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(328, 4).Value = 19
$ws1.Cells.Item(350, 1).Value = '20:12:03'
$ws1.Cells.Item(350, 2).Value = '20:57'
$ws1.Cells.Item(350, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(350, 4).Value = 45
$ws1.Cells.Item(351, 1).Value = '19:11:45'
$ws1.Cells.Item(351, 2).Value = '21:06'
$ws1.Cells.Item(351, 4).Value = 115
$ws1.Cells.Item(352, 1).Value = '19:48:11'
$ws1.Cells.Item(352, 2).Value = '21:07'
$ws1.Cells.Item(352, 3).Value = '10_OLMOS'
$ws1.Cells.Item(352, 4).Value = 79
$ws1.Cells.Item(353, 1).Value = '19:35:56'
$ws1.Cells.Item(353, 2).Value = '21:09'
$ws1.Cells.Item(353, 4).Value = 94
$ws1.Cells.Item(354, 1).Value = '19:48:11'
$ws1.Cells.Item(354, 2).Value = '21:10'
$ws1.Cells.Item(354, 3).Value = '15_ABASTO'
$ws1.Cells.Item(354, 4).Value = 82
$ws1.Cells.Item(355, 2).Value = '21:28'
$ws1.Cells.Item(355, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(355, 4).Value = 113
$ws1.Cells.Item(356, 1).Value = '19:35:56'
$ws1.Cells.Item(356, 2).Value = '21:33'
$ws1.Cells.Item(356, 4).Value = 118
$ws1.Cells.Item(357, 1).Value = '19:48:11'
$ws1.Cells.Item(357, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(357, 4).Value = 106
$ws1.Cells.Item(358, 1).Value = '19:55:23'
$ws1.Cells.Item(358, 2).Value = '21:34'
$ws1.Cells.Item(358, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(358, 4).Value = 99
$ws1.Cells.Item(359, 1).Value = '20:12:03'
$ws1.Cells.Item(359, 2).Value = '21:38'
$ws1.Cells.Item(359, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(359, 4).Value = 86
$ws1.Cells.Item(359, 5).Value = 'LP1912'
$ws1.Cells.Item(360, 1).Value = '19:48:11'
$ws1.Cells.Item(360, 2).Value = '21:46'
$ws1.Cells.Item(360, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(360, 4).Value = 118
$ws1.Cells.Item(360, 5).Value = 'LP1912'
$ws1.Cells.Item(361, 1).Value = '20:12:03'
$ws1.Cells.Item(361, 2).Value = '21:47'
$ws1.Cells.Item(361, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(361, 4).Value = 95
$ws1.Cells.Item(361, 5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 20:12:03'

# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 20:12:03'
$ws3.Cells.Item(3, 1).Value = 'Total filas: 44'
$ws3.Cells.Item(49, 1).Value = '20:12:03'
$ws3.Cells.Item(49, 2).Value = '21:30'
$ws3.Cells.Item(49, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(49, 4).Value = 78
$ws3.Cells.Item(49, 5).Value = 'L6203'
